$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$dst = $ws.Range("S12")
$dst.Borders(7).LineStyle = 1
$dst.Borders(10).LineStyle = 1
$dst.Borders(8).LineStyle = 1
$dst.Borders(9).LineStyle = 1
$dst.Borders(8).LineStyle = -4142
$dst.Borders(9).LineStyle = -4142
Write-Host "done"
